$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F9").Value = 'In column ''CG'' of the extraction file, variable ID should be ''RC-7'' when the actual value is ''R-C-7'''
$ws.Range("F10").Value = 'In column ''BM'' of the extraction file, variable label should be ''p-value'' when the actual value is ''pvalue'''
$ws.Range("F14").Value = 'In column ''BW'' of the extraction file, variable ID should be ''RI-5'' when the actual value is ''RI-58'''
$ws.Range("F15").Value = 'In column ''CZ'' of the extraction file, variable label should be ''Timepoint unit'' when the actual value is ''Invalid unit'''
$ws.Range("F18").Value = 'In column ''DM'' of the extraction file, variable ID should be ''ST-9'' when the actual value is ''9-ST'''
$ws.Range("F19").Value = 'In column ''BT'' of the extraction file, variable label should be ''Timepoint'' when the actual value is ''Timepointssss'''
$ws.Range("F20").Value = 'In column ''BD'' of the extraction file, variable label should be ''Number of events'' when the actual value is ''% of patients experiencing event'''
$ws.Range("F24").Value = 'In column ''DR'' of the extraction file, variable ID should be ''T-2'' when the actual value is ''T-22'''
$ws.Range("F25").Value = 'In column ''AT'' of the extraction file, variable label should be ''Add Endpoint Description custom columns'' when the actual value is ''Invalid Col Name'''
$ws.Range("F32").Value = 'In column ''BZ'' of the extraction file, variable label should be ''Add Endpoint overall results custom columns'' when the actual value is ''Add Endpoint overall results custom colIvalidumns'''
$ws.Range("F51").Value = 'Column C-4 mapping key ''Ext_CustomData'' is not supported. Only "Arm_" is supported for arm-related data or "Extraction_" is supported for extraction-related data. Please modify accordingly and re-upload.'

$ws.Range("F25").Select()
